# Apply "想去人数" (F column) count updates across the four sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1870
$ws1.Range("F5").Value  = 871
$ws1.Range("F6").Value  = 392
$ws1.Range("F7").Value  = 759
$ws1.Range("F8").Value  = 13272
$ws1.Range("F9").Value  = 13162
$ws1.Range("F10").Value = 1011
$ws1.Range("F12").Value = 19
$ws1.Range("F13").Value = 554
$ws1.Range("F15").Value = 665
$ws1.Range("F19").Value = 63
$ws1.Range("F21").Value = 234

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 33
$ws2.Range("F9").Value = 24

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 32

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1870
$ws4.Range("F6").Value  = 871
$ws4.Range("F7").Value  = 392
$ws4.Range("F9").Value  = 759
$ws4.Range("F10").Value = 13272
$ws4.Range("F11").Value = 13162
$ws4.Range("F12").Value = 1011
$ws4.Range("F14").Value = 19
$ws4.Range("F15").Value = 554
$ws4.Range("F17").Value = 665
$ws4.Range("F18").Value = 33
$ws4.Range("F23").Value = 63
$ws4.Range("F27").Value = 32
$ws4.Range("F28").Value = 234
$ws4.Range("F34").Value = 24
